$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values
# (e.g. "0.9998", "254.49", "1.000") are not auto-converted to numbers,
# matching the inlineStr/text representation used in the workbook.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '31.247.09'
$ws.Range('E2').Value = '  +2.13%  '
$ws.Range('D3').Value = '1.997.47'
$ws.Range('E3').Value = '  +6.02%  '
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '0.7755'
$ws.Range('E5').Value = '  +63.76%  '
$ws.Range('D6').Value = '254.49'
$ws.Range('E6').Value = '  +3.33%  '
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '0.3481'
$ws.Range('E8').Value = '  +20.26%  '
$ws.Range('D9').Value = '27.89'
$ws.Range('E9').Value = '  +24.71%  '
$ws.Range('E10').Value = '  +7.98%  '
$ws.Range('D11').Value = '0.8411'
$ws.Range('E11').Value = '  +10.27%  '
$ws.Range('D12').Value = '0.08203'
$ws.Range('E12').Value = '  +4.86%  '
$ws.Range('B13').Value = 'Litecoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D13').Value = '100.75'
$ws.Range('E13').Value = '  +0.98%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.995.74'
$ws.Range('E14').Value = '  +5.98%  '
$ws.Range('E15').Value = '  +7.47%  '
$ws.Range('E16').Value = '  +15.37%  '
$ws.Range('D17').Value = '272.73'
$ws.Range('E17').Value = '  -3.97%  '
$ws.Range('D18').Value = '31.244.81'
$ws.Range('E18').Value = '  +2.20%  '
$ws.Range('D19').Value = '5.972'
$ws.Range('E19').Value = '  +11.64%  '
$ws.Range('D20').Value = '0.000008008'
$ws.Range('E20').Value = '  +6.49%  '
$ws.Range('D21').Value = '2.258.20'
$ws.Range('D22').Value = '1.0000'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').Value = '0.9992'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').Value = '7.094'
$ws.Range('E24').Value = '  +10.30%  '
$ws.Range('D25').Value = '10.00'
$ws.Range('E25').Value = '  +9.06%  '
$ws.Range('D26').Value = '165.27'
$ws.Range('E26').Value = '  +0.93%  '
$ws.Range('D27').Value = '0.1423'
$ws.Range('E27').Value = '  +45.80%  '
$ws.Range('D28').Value = '19.90'
$ws.Range('E28').Value = '  +4.69%  '
$ws.Range('D29').Value = '2.402'
$ws.Range('E29').Value = '  +25.98%  '
$ws.Range('D30').Value = '1.596'
$ws.Range('E30').Value = '  +6.34%  '
$ws.Range('D31').Value = '1.365'
$ws.Range('E31').Value = '  +2.82%  '
$ws.Range('D32').Value = '4.604'
$ws.Range('E32').Value = '  +8.25%  '
$ws.Range('E33').Value = '  +6.04%  '
$ws.Range('D34').Value = '0.05326'
$ws.Range('E34').Value = '  +9.97%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.246'
$ws.Range('E35').Value = '  +10.19%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.7922'
$ws.Range('E36').Value = '  +13.28%  '
$ws.Range('D37').Value = '2.769'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D38').Value = '0.9992'
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.02005'
$ws.Range('E39').Value = '  +5.23%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.914'
$ws.Range('E40').Value = '  +1.61%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '82.51'
$ws.Range('E41').Value = '  +9.17%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '6.784'
$ws.Range('E42').Value = '  +7.40%  '
$ws.Range('D43').Value = '0.4677'
$ws.Range('E43').Value = '  +9.97%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '2.133'
$ws.Range('E44').Value = '  +7.81%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = '0.8564'
$ws.Range('E45').Value = '  +2.15%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = '105.21'
$ws.Range('E46').Value = '  +3.66%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').Value = '1.001'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '10.04'
$ws.Range('E48').Value = '  +0.51%  '
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').Value = '7.720'
$ws.Range('E49').Value = '  +10.03%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = '37.49'
$ws.Range('E50').Value = '  +6.26%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '1.564'
$ws.Range('E51').Value = '  +16.65%  '

# Restore default cell style on column D so no stray style/number-format
# metadata is introduced (the values are already stored as text strings).
$dRange.Style = "Normal"
